$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 5
$ws.Range('C2').Value = 72.791
$ws.Range('D2').Value = 'MASTER MEATS (BOUTIQUE DE CARNES)'
$ws.Range('E2').Value = 'Rua Abílio Soares, 731 - Paraíso São Paulo/SP CEP:04005003'
$ws.Range('F2').Value = 3

$ws.Range('A3').Value = 11
$ws.Range('C3').Value = 72.838
$ws.Range('D3').Value = 'MYW1 O BAR LTDA (VASSOURA QUEBRADA - PERDIZES)'
$ws.Range('E3').Value = 'Rua Desembargador do Vale, 836, ANEXO 830 - Perdizes São Paulo/SP CEP:05010040'
$ws.Range('F3').Value = 10

$ws.Range('A4').Value = 14
$ws.Range('C4').Value = 72.859
$ws.Range('D4').Value = 'MRL BUS COM. DE ALIM. EIRELI - ME (BUSGER - VILA MADALENA)'
$ws.Range('E4').Value = 'Rua Alves Guimarães, 1091, COZINHA 4 - Pinheiros São Paulo/SP CEP:05410-002'
$ws.Range('F4').Value = 15

$ws.Range('A5').Value = 10
$ws.Range('C5').Value = 72.891
$ws.Range('D5').Value = 'ECULLY CHARBON RESTAURANTE LTDA (ECULLY CHARBON)'
$ws.Range('E5').Value = 'Rua Doutor Augusto de Miranda, 549 - Vila Pompéia São Paulo/SP CEP:05026000'
$ws.Range('F5').Value = 14

$ws.Range('A6').Value = 7
$ws.Range('C6').Value = 72.903
$ws.Range('D6').Value = 'BAR & REST PICANHARIA DOS AMIGOS LTDA (PICANHARIA DOS AMIGOS - VILA LEOPOLDINA)'
$ws.Range('E6').Value = 'Rua Guaipá, 1017,  - Vila Leopoldina São Paulo/SP CEP:05089-001'
$ws.Range('F6').Value = 2

$ws.Range('A7').Value = 9
$ws.Range('C7').Value = 72.935
$ws.Range('D7').Value = 'CAPITAO COM. E DIST. BEBIDAS E ALIMENTOS (CAPITAO BARLEY)'
$ws.Range('E7').Value = 'Rua Coriolano, 301 - Vila Romana São Paulo/SP CEP:05047001'
$ws.Range('F7').Value = 2

$ws.Range('A8').Value = 6
$ws.Range('C8').Value = 72.937
$ws.Range('D8').Value = 'PARCEL SW BURGUER LTDA (N! BURGER - LAPA)'
$ws.Range('E8').Value = 'Rua Catão, 479, NBURGER - Vila Romana São Paulo/SP CEP:05049000'
$ws.Range('F8').Value = 11

$ws.Range('A9').Value = 8
$ws.Range('C9').Value = 72.941
$ws.Range('D9').Value = 'TOSQUINHO LANCHES LTDA (TOSQUINHO LANCHES)'
$ws.Range('E9').Value = 'RUA CAMILO, 763, sem complemento - VILA ROMANA São Paulo/SP CEP:05045020'
$ws.Range('F9').Value = 8

$ws.Range('A10').Value = 1
$ws.Range('C10').Value = 72.956
$ws.Range('D10').Value = '*CLIENTE AMOSTRA (CLIENTE AMOSTRA)*'
$ws.Range('E10').Value = 'Rua José Mariano Filho, 200,  - Jardim Oriental São Paulo/SP CEP:04347-180'
$ws.Range('F10').Value = 4

$ws.Range('A11').Value = 4
$ws.Range('C11').Value = 72.967
$ws.Range('D11').Value = 'GILBERTO CAMPOS DE AZAMBUJA ME (ROYAL MEAT - PARAISO)'
$ws.Range('E11').Value = 'Rua Doutor Tomás Carvalhal, 626 - Paraíso São Paulo/SP CEP:04006001'
$ws.Range('F11').Value = 5

$ws.Range('A12').Value = 12
$ws.Range('C12').Value = 72.97
$ws.Range('D12').Value = 'THE BEAR BURGER REST. LTDA EPP (THE BEAR BURGER)'
$ws.Range('E12').Value = 'Rua Caraíbas, 964, IMOBILIARIA ESTEVAM - Perdizes São Paulo/SP CEP:05020000'
$ws.Range('F12').Value = 5

$ws.Range('A13').Value = 3
$ws.Range('C13').Value = 72.988
$ws.Range('D13').Value = 'BUSGER COM. DE ALIM. LTDA (BUSGER - KLABIN)'
$ws.Range('E13').Value = 'Rua Vergueiro, 4289,  - Vila Mariana São Paulo/SP CEP:04101-901'
$ws.Range('F13').Value = 30

$ws.Range('A14').Value = 2
$ws.Range('C14').Value = 72.99
$ws.Range('D14').Value = 'ESTEFOODS COM. DE ALIMENTOS LTDA (BUSGER - BORGES LAGOA)'
$ws.Range('E14').Value = 'Rua Borges Lagoa, 1050,  - Vila Clementino São Paulo/SP CEP:04038-002'
$ws.Range('F14').Value = 26

$ws.Range('A15').Value = 13
$ws.Range('C15').Value = 73.008
$ws.Range('D15').Value = 'BORGER BURGER LTDA (BORGER - PERDIZES)'
$ws.Range('E15').Value = 'Rua Cardoso de Almeida, 587,  - Perdizes São Paulo/SP CEP:05013-000'
$ws.Range('F15').Value = 2

